$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.753.95"
$ws.Range("E2").Value = "  -2.54%  "
$ws.Range("D3").Value = "2.445.69"
$ws.Range("E3").Value = "  -3.84%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'522.94"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.71%  "
$ws.Range("D6").Value = "'129.03"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.42%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").Value = "'0.562"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.10%  "
$ws.Range("D9").Value = "'0.0973"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.82%  "
$ws.Range("E10").Value = "  -2.24%  "
$ws.Range("E11").Value = "  -5.51%  "
$ws.Range("E12").Value = "  -4.48%  "
$ws.Range("D13").Value = "2.877.75"
$ws.Range("E13").Value = "  -3.94%  "
$ws.Range("D14").Value = "57.681.48"
$ws.Range("E14").Value = "  -2.56%  "
$ws.Range("E15").Value = "  -4.08%  "
$ws.Range("E16").Value = "  -3.57%  "
$ws.Range("D17").Value = "2.445.27"
$ws.Range("E17").Value = "  -3.85%  "
$ws.Range("D18").Value = "'10.35"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -3.78%  "
$ws.Range("E19").Value = "  -2.60%  "
$ws.Range("D20").Value = "'311.74"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.74%  "
$ws.Range("D21").Value = "'6.08"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.25%  "
$ws.Range("E22").Value = "  -0.11%  "
$ws.Range("D23").Value = "'64.92"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.63%  "
$ws.Range("E24").Value = "  -2.63%  "
$ws.Range("D26").Value = "2.565.17"
$ws.Range("E26").Value = "  -3.72%  "
$ws.Range("D28").Value = "'7.22"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.17%  "
$ws.Range("D29").Value = "'173.18"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.04%  "
$ws.Range("E30").Value = "  -3.53%  "
$ws.Range("D32").Value = "'6.13"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.72%  "
$ws.Range("D33").Value = "'1.12"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -9.00%  "
$ws.Range("E35").Value = "  -0.15%  "
$ws.Range("D36").Value = "'17.78"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.86%  "
$ws.Range("E37").Value = "  -7.30%  "
$ws.Range("D38").Value = "'3.75"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -5.85%  "
$ws.Range("D39").Value = "'36.27"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.38%  "
$ws.Range("D40").Value = "'0.800"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.38%  "
$ws.Range("E41").Value = "  -5.34%  "
$ws.Range("D42").Value = "'3.37"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.49%  "
$ws.Range("E43").Value = "  -3.68%  "
$ws.Range("D44").Value = "'4.76"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -6.97%  "
$ws.Range("B45").Value = "Stellar"
$ws.Range("C45").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D45").Value = "'0.0917"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.80%  "
$ws.Range("B46").Value = "Bittensor"
$ws.Range("C46").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D46").Value = "'252.59"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -9.85%  "
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").Value = "'121.19"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -9.61%  "
$ws.Range("E48").Value = "  -3.70%  "
$ws.Range("D50").Value = "'16.95"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -5.47%  "
$ws.Range("D51").Value = "'16.09"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -6.45%  "
